$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) updates to column F ("想去人数")
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 3846
$ws.Range("F6").Value = 3842
$ws.Range("F7").Value = 392
$ws.Range("F8").Value = 208
$ws.Range("F9").Value = 60
$ws.Range("F10").Value = 8777
$ws.Range("F11").Value = 83
$ws.Range("F12").Value = 137
$ws.Range("F13").Value = 308
$ws.Range("F14").Value = 347
$ws.Range("F16").Value = 108
$ws.Range("F17").Value = 9
$ws.Range("F18").Value = 378
$ws.Range("F19").Value = 11123
$ws.Range("F21").Value = 294
$ws.Range("F28").Value = 195
$ws.Range("F34").Value = 2084
$ws.Range("F38").Value = 912
$ws.Range("F39").Value = 2571
$ws.Range("F40").Value = 288
$ws.Range("F42").Value = 1255
$ws.Range("F44").Value = 758
$ws.Range("F46").Value = 355
$ws.Range("F48").Value = 92
$ws.Range("F49").Value = 86

# Sheet "演出" (sheet2) updates to column F ("想去人数")
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F10").Value = 51
$ws.Range("F24").Value = 27

# Sheet "全部类型" (sheet4) updates to column F ("想去人数")
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 3846
$ws.Range("F6").Value = 3842
$ws.Range("F7").Value = 392
$ws.Range("F9").Value = 51
$ws.Range("F10").Value = 208
$ws.Range("F11").Value = 60
$ws.Range("F12").Value = 8778
$ws.Range("F14").Value = 83
$ws.Range("F15").Value = 137
$ws.Range("F16").Value = 308
$ws.Range("F17").Value = 347
$ws.Range("F19").Value = 108
$ws.Range("F20").Value = 378
$ws.Range("F21").Value = 11123
$ws.Range("F22").Value = 294
$ws.Range("F26").Value = 195
$ws.Range("F30").Value = 2084
$ws.Range("F34").Value = 912
$ws.Range("F37").Value = 2571
$ws.Range("F38").Value = 288
$ws.Range("F41").Value = 1255
$ws.Range("F43").Value = 758
$ws.Range("F45").Value = 355
$ws.Range("F46").Value = 27
$ws.Range("F48").Value = 92
$ws.Range("F49").Value = 86
